# Refresh the scraped "Price" (D) / "Volume(1h)" (E) columns of the cryptos
# worksheet with the latest values from the GitHub Actions run.
#
# Both columns hold plain text (not numbers): the Price column mixes
# thousands-dot-grouped values ("29.381.51") with plain decimals ("242.46"),
# and the Volume column is a padded, signed percentage string
# ("  +0.07%  "). Values that could otherwise be auto-parsed by Excel as a
# plain number are written with a leading apostrophe so they stick as text,
# matching the original inline-string cell contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.381.51'

# Row 3
$ws.Range("D3").Value = '1.881.67'
$ws.Range("E3").Value = '  +0.29%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").Value = '''0.7129'
$ws.Range("E5").Value = '  -0.05%  '

# Row 6
$ws.Range("D6").Value = '''242.46'
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").Value = '''0.08057'
$ws.Range("E8").Value = '  +3.94%  '

# Row 9
$ws.Range("D9").Value = '''0.3129'
$ws.Range("E9").Value = '  +0.64%  '

# Row 10
$ws.Range("D10").Value = '''25.26'
$ws.Range("E10").Value = '  +1.27%  '

# Row 11
$ws.Range("D11").Value = '''0.08327'
$ws.Range("E11").Value = '  -2.41%  '

# Row 12
$ws.Range("D12").Value = '1.893.11'
$ws.Range("E12").Value = '  +0.46%  '

# Row 13
$ws.Range("D13").Value = '''5.247'
$ws.Range("E13").Value = '  +0.61%  '

# Row 14
$ws.Range("D14").Value = '''0.7189'
$ws.Range("E14").Value = '  +1.19%  '

# Row 15
$ws.Range("D15").Value = '''93.85'
$ws.Range("E15").Value = '  +2.62%  '

# Row 16
$ws.Range("D16").Value = '''6.332'
$ws.Range("E16").Value = '  +5.32%  '

# Row 17
$ws.Range("D17").Value = '''0.000008560'
$ws.Range("E17").Value = '  +4.53%  '

# Row 18
$ws.Range("D18").Value = '29.392.36'
$ws.Range("E18").Value = '  +0.29%  '

# Row 19
$ws.Range("D19").Value = '''242.15'
$ws.Range("E19").Value = '  -0.10%  '

# Row 20
$ws.Range("D20").Value = '2.138.24'
$ws.Range("E20").Value = '  -0.11%  '

# Row 21
$ws.Range("D21").Value = '''13.25'
$ws.Range("E21").Value = '  -0.02%  '

# Row 22
$ws.Range("E22").Value = '  +0.11%  '

# Row 23
$ws.Range("E23").Value = '  +0.51%  '

# Row 24
$ws.Range("E24").Value = '  +0.06%  '

# Row 25
$ws.Range("D25").Value = '''0.1592'
$ws.Range("E25").Value = '  -0.88%  '

# Row 26
$ws.Range("D26").Value = '''163.73'
$ws.Range("E26").Value = '  +0.47%  '

# Row 27
$ws.Range("D27").Value = '''9.075'
$ws.Range("E27").Value = '  +0.29%  '

# Row 28
$ws.Range("E28").Value = '  +0.48%  '

# Row 29
$ws.Range("D29").Value = '''1.509'
$ws.Range("E29").Value = '  -0.26%  '

# Row 30
$ws.Range("D30").Value = '''4.413'
$ws.Range("E30").Value = '  +0.27%  '

# Row 31
$ws.Range("D31").Value = '''4.337'
$ws.Range("E31").Value = '  +0.43%  '

# Row 32
$ws.Range("D32").Value = '''1.201'
$ws.Range("E32").Value = '  -6.17%  '

# Row 33
$ws.Range("D33").Value = '''0.05379'
$ws.Range("E33").Value = '  +2.48%  '

# Row 34
$ws.Range("D34").Value = '''1.948'
$ws.Range("E34").Value = '  +0.79%  '

# Row 35
$ws.Range("D35").Value = '''1.182'
$ws.Range("E35").Value = '  +0.59%  '

# Row 36
$ws.Range("D36").Value = '''0.7487'
$ws.Range("E36").Value = '  +0.21%  '

# Row 37
$ws.Range("D37").Value = '''2.699'
$ws.Range("E37").Value = '  +0.47%  '

# Row 38
$ws.Range("D38").Value = '''0.01888'
$ws.Range("E38").Value = '  +1.10%  '

# Row 39
$ws.Range("D39").Value = '1.287.47'
$ws.Range("E39").Value = '  +8.97%  '

# Row 40
$ws.Range("D40").Value = '''2.746'
$ws.Range("E40").Value = '  +1.11%  '

# Row 41
$ws.Range("D41").Value = '''6.593'
$ws.Range("E41").Value = '  +3.18%  '

# Row 42
$ws.Range("D42").Value = '''0.9191'
$ws.Range("E42").Value = '  +3.63%  '

# Row 43
$ws.Range("D43").Value = '''74.90'
$ws.Range("E43").Value = '  +2.67%  '

# Row 44
$ws.Range("D44").Value = '''111.66'
$ws.Range("E44").Value = '  +4.83%  '

# Row 45
$ws.Range("E45").Value = '  +0.09%  '

# Row 46
$ws.Range("E46").Value = '  +5.43%  '

# Row 47
$ws.Range("D47").Value = '2.032.57'
$ws.Range("E47").Value = '  +0.04%  '

# Row 48
$ws.Range("E48").Value = '  -0.10%  '

# Row 49
$ws.Range("D49").Value = '''0.5222'
$ws.Range("E49").Value = '  +0.29%  '

# Row 50
$ws.Range("D50").Value = '''9.526'
$ws.Range("E50").Value = '  +1.46%  '

# Row 51
$ws.Range("D51").Value = '''0.4390'
$ws.Range("E51").Value = '  +1.78%  '
